$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot in columns D and E.
#
# Column D holds price text that often LOOKS like a plain number (e.g.
# "231.68"), while other rows use a dotted, European-style thousands
# grouping that is not a valid number (e.g. "34.862.29"). In the source
# workbook every one of these cells is stored as text, not a number.
# Simply assigning a numeric-looking string to `.Value` lets Excel
# auto-convert it to a real number, which would change the cell's
# underlying type. To keep these cells as text (matching the original
# file), force a text number format before writing the value, then put
# the style back to "Normal" so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.862.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.809.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.19%  '
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.313'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0682'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.54%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.070.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.801.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.67'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.655'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.821.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0783'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.34%  '
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '173.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +31.83%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.339.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0549'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.15%  '
$ws.Range("E34").Value = '  +1.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '93.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.675'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0192'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.305.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("E41").Value = '  +4.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.984'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.20%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0512'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.986.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("E51").Value = '  +5.76%  '
